$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "L/min" unit label with the degree symbol "°" in B2 and D2
$ws.Range("B2").Value = "°"
$ws.Range("D2").Value = "°"

# Move the selection to D2 (matches the final sheetView selection in the diff)
$ws.Range("D2").Select()
